$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.333402
$ws.Range("H2").Value = 73.000206
$ws.Range("I2").Value = 0.9697175080062574
$ws.Range("J2").Value = 0.9697175080062576
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 117.956133750464
$ws.Range("R2").Value = 1061.605203754176
$ws.Range("S2").Value = 0.03297690881271426
$ws.Range("T2").Value = 0.03297690881271426

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.333402
$ws.Range("H3").Value = 73.000206
$ws.Range("I3").Value = 0.9697175080062574
$ws.Range("J3").Value = 0.9697175080062576
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 39.59271617152001
$ws.Range("R3").Value = 356.33444554368
$ws.Range("S3").Value = 0.01106890629018057
$ws.Range("T3").Value = 0.01106890629018058

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.333402
$ws.Range("H4").Value = 73.000206
$ws.Range("I4").Value = 0.9697175080062574
$ws.Range("J4").Value = 0.9697175080062576
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 3311.063951968409
$ws.Range("R4").Value = 29799.57556771568
$ws.Range("S4").Value = 0.9256716929033626
$ws.Range("T4").Value = 0.9256716929033627

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.288887
$ws.Range("H5").Value = 0.866661
$ws.Range("I5").Value = 0.01151252018667195
$ws.Range("J5").Value = 0.01151252018667195
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 1.400379347317333
$ws.Range("R5").Value = 12.603414125856
$ws.Range("S5").Value = 0.0003915030153275971
$ws.Range("T5").Value = 0.0003915030153275971

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.288887
$ws.Range("H6").Value = 0.866661
$ws.Range("I6").Value = 0.01151252018667195
$ws.Range("J6").Value = 0.01151252018667195
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 0.4700461117866667
$ws.Range("R6").Value = 4.23041500608
$ws.Range("S6").Value = 0.0001314104427918215
$ws.Range("T6").Value = 0.0001314104427918215

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.288887
$ws.Range("H7").Value = 0.866661
$ws.Range("I7").Value = 0.01151252018667195
$ws.Range("J7").Value = 0.01151252018667195
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 39.30906709601467
$ws.Range("R7").Value = 353.781603864132
$ws.Range("S7").Value = 0.01098960672855253
$ws.Range("T7").Value = 0.01098960672855253

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4710003333333333
$ws.Range("H8").Value = 1.413001
$ws.Range("I8").Value = 0.01876997180707065
$ws.Range("J8").Value = 0.01876997180707065
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.847498666666667
$ws.Range("N8").Value = 14.542496
$ws.Range("O8").Value = 0.03400671694637637
$ws.Range("P8").Value = 0.03400671694637637
$ws.Range("Q8").Value = 2.283173487832889
$ws.Range("R8").Value = 20.548561390496
$ws.Range("S8").Value = 0.0006383051183345159
$ws.Range("T8").Value = 0.000638305118334516

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4710003333333333
$ws.Range("H9").Value = 1.413001
$ws.Range("I9").Value = 0.01876997180707065
$ws.Range("J9").Value = 0.01876997180707065
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.627093333333334
$ws.Range("N9").Value = 4.88128
$ws.Range("O9").Value = 0.01141456784970118
$ws.Range("P9").Value = 0.01141456784970118
$ws.Range("Q9").Value = 0.7663615023644444
$ws.Range("R9").Value = 6.897253521280001
$ws.Range("S9").Value = 0.0002142511167287862
$ws.Range("T9").Value = 0.0002142511167287862

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4710003333333333
$ws.Range("H10").Value = 1.413001
$ws.Range("I10").Value = 0.01876997180707065
$ws.Range("J10").Value = 0.01876997180707065
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.0707373333333
$ws.Range("N10").Value = 408.212212
$ws.Range("O10").Value = 0.9545787152039225
$ws.Range("P10").Value = 0.9545787152039225
$ws.Range("Q10").Value = 64.08936264091244
$ws.Range("R10").Value = 576.804263768212
$ws.Range("S10").Value = 0.01791741557200734
$ws.Range("T10").Value = 0.01791741557200735
